$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 195
$ws.Range("I4").Value = 195
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 195
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -81
$ws.Range("N4").ClearContents()

$ws.Range("H51").Value = 1867.25
$ws.Range("I51").Value = 1989.6666
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 1989.6666
$ws.Range("L51").Value = 1500
$ws.Range("M51").Value = -1505.6666
$ws.Range("N51").Value = -2468

$ws.Range("H74").Value = 3868.5625
$ws.Range("I74").Value = 3742.7144
$ws.Range("J74").Value = 3966.4443
$ws.Range("K74").Value = 3742.7144
$ws.Range("L74").Value = 3966.4443
$ws.Range("M74").Value = -2806.7144
$ws.Range("N74").Value = -5838.4443

$ws.Range("H77").Value = 3868.5625
$ws.Range("I77").Value = 3742.7144
$ws.Range("J77").Value = 3966.4443
$ws.Range("K77").Value = 18713.572
$ws.Range("L77").Value = 19832.2215
$ws.Range("M77").Value = -14033.572
$ws.Range("N77").Value = -29192.2215

$ws.Range("H113").Value = 123133.89
$ws.Range("I113").Value = 153457.86
$ws.Range("J113").Value = 17000
$ws.Range("K113").Value = 153457.86
$ws.Range("L113").Value = 17000
$ws.Range("M113").Value = -150203.86
$ws.Range("N113").Value = -23508

$ws.Range("H125").Value = 18685702
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 22422742
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 201804678
$ws.Range("M125").Value = -2040
$ws.Range("N125").Value = -201809598

$ws.Range("H132").Value = 374610.7
$ws.Range("I132").Value = 487126.28
$ws.Range("J132").Value = 62067.332
$ws.Range("K132").Value = 1461378.84
$ws.Range("L132").Value = 186201.996
$ws.Range("M132").Value = -1458848.84
$ws.Range("N132").Value = -191261.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1639.9166
$ws.Range("I2").Value = 843.8
$ws.Range("J2").Value = 2208.5715
$ws.Range("K2").Value = 843.8
$ws.Range("L2").Value = 2208.5715
$ws.Range("M2").Value = -730.8
$ws.Range("N2").Value = -2434.5715

$ws.Range("H4").Value = 236.6
$ws.Range("I4").Value = 236.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 236.6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -120.6
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 13725.829
$ws.Range("I32").Value = 1297.3513
$ws.Range("K32").Value = 1297.3513
$ws.Range("M32").Value = -1010.3513

$ws.Range("H116").Value = 1639.9166
$ws.Range("I116").Value = 843.8
$ws.Range("J116").Value = 2208.5715
$ws.Range("K116").Value = 843.8
$ws.Range("L116").Value = 2208.5715
$ws.Range("M116").Value = 1450.2
$ws.Range("N116").Value = -6796.5715

$ws.Range("H132").Value = 3003.625
$ws.Range("I132").Value = 2737.077
$ws.Range("K132").Value = 8211.231
$ws.Range("M132").Value = -5681.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1639.9166
$ws.Range("I3").Value = 843.8
$ws.Range("J3").Value = 2208.5715
$ws.Range("K3").Value = 843.8
$ws.Range("L3").Value = 2208.5715
$ws.Range("M3").Value = -729.8
$ws.Range("N3").Value = -2436.5715

$ws.Range("H107").Value = 885.3684
$ws.Range("I107").Value = 686.1429000000001
$ws.Range("K107").Value = 686.1429000000001
$ws.Range("M107").Value = 1233.8571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 227000.8
$ws.Range("I13").Value = 100004
$ws.Range("J13").Value = 258750
$ws.Range("K13").Value = 100004
$ws.Range("L13").Value = 258750
$ws.Range("M13").Value = -99865
$ws.Range("N13").Value = -259028

$ws.Range("H99").Value = 5683152
$ws.Range("I99").Value = 7813615.5
$ws.Range("J99").Value = 1916.6666
$ws.Range("K99").Value = 7813615.5
$ws.Range("L99").Value = 1916.6666
$ws.Range("M99").Value = -7812117.5
$ws.Range("N99").Value = -4912.6666

$ws.Range("H107").Value = 206.38889
$ws.Range("I107").Value = 94.333336
$ws.Range("J107").Value = 430.5
$ws.Range("K107").Value = 94.333336
$ws.Range("L107").Value = 430.5
$ws.Range("M107").Value = 1825.666664
$ws.Range("N107").Value = -4270.5

$ws.Range("H126").Value = 5683152
$ws.Range("I126").Value = 7813615.5
$ws.Range("J126").Value = 1916.6666
$ws.Range("K126").Value = 23440846.5
$ws.Range("L126").Value = 5749.9998
$ws.Range("M126").Value = -23438376.5
$ws.Range("N126").Value = -10689.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 139.125
$ws.Range("I2").Value = 26.076923
$ws.Range("J2").Value = 272.72726
$ws.Range("K2").Value = 156.461538
$ws.Range("L2").Value = 1636.36356
$ws.Range("M2").Value = -43.46153800000002
$ws.Range("N2").Value = -1862.36356

$ws.Range("H33").Value = 426.2857
$ws.Range("I33").Value = 95.5
$ws.Range("J33").Value = 558.6
$ws.Range("K33").Value = 573
$ws.Range("L33").Value = 3351.6
$ws.Range("M33").Value = -290
$ws.Range("N33").Value = -3917.6

$ws.Range("H80").Value = 1108.8889
$ws.Range("J80").Value = 1108.8889
$ws.Range("L80").Value = 3326.6667
$ws.Range("N80").Value = -5198.6667

$ws.Range("H82").Value = 2526.5
$ws.Range("I82").Value = 404
$ws.Range("J82").Value = 3800
$ws.Range("K82").Value = 1212
$ws.Range("L82").Value = 11400
$ws.Range("M82").Value = -806
$ws.Range("N82").Value = -12212

$ws.Range("H83").Value = 1108.8889
$ws.Range("J83").Value = 1108.8889
$ws.Range("L83").Value = 9980.000099999999
$ws.Range("N83").Value = -19340.0001

$ws.Range("H85").Value = 2526.5
$ws.Range("I85").Value = 404
$ws.Range("J85").Value = 3800
$ws.Range("K85").Value = 1212
$ws.Range("L85").Value = 11400
$ws.Range("M85").Value = 192
$ws.Range("N85").Value = -14208

$ws.Range("H113").Value = 7247027.5
$ws.Range("I113").Value = 494.53845
$ws.Range("J113").Value = 11628652
$ws.Range("K113").Value = 1483.61535
$ws.Range("L113").Value = 34885956
$ws.Range("M113").Value = 686.38465
$ws.Range("N113").Value = -34890296

$ws.Range("H122").Value = 783.7692
$ws.Range("I122").Value = 303
$ws.Range("J122").Value = 997.44446
$ws.Range("K122").Value = 2727
$ws.Range("L122").Value = 8977.00014
$ws.Range("M122").Value = -277
$ws.Range("N122").Value = -13877.00014

$ws.Range("H132").Value = 1094.7778
$ws.Range("I132").Value = 779.63635
$ws.Range("J132").Value = 1311.4375
$ws.Range("K132").Value = 7016.72715
$ws.Range("L132").Value = 11802.9375
$ws.Range("M132").Value = -4486.72715
$ws.Range("N132").Value = -16862.9375

$ws.Range("H133").Value = 9575
$ws.Range("I133").Value = 2300
$ws.Range("K133").Value = 6900
$ws.Range("M133").Value = -1840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38.866665
$ws.Range("I2").Value = 31.3
$ws.Range("J2").Value = 54
$ws.Range("K2").Value = 31.3
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 81.7
$ws.Range("N2").Value = -280

$ws.Range("H80").Value = 2868.6667
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 2953
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 2953
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4949

$ws.Range("H83").Value = 2868.6667
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 2953
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 14765
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -24749

$ws.Range("H123").Value = 9184.708000000001
$ws.Range("J123").Value = 9184.708000000001
$ws.Range("L123").Value = 9184.708000000001
$ws.Range("N123").Value = -14084.708

$ws.Range("H132").Value = 3271.5386
$ws.Range("I132").Value = 2941.4443
$ws.Range("J132").Value = 3554.476
$ws.Range("K132").Value = 8824.332900000001
$ws.Range("L132").Value = 10663.428
$ws.Range("M132").Value = -6294.332900000001
$ws.Range("N132").Value = -15723.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3044.1177
$ws.Range("I40").Value = 790
$ws.Range("J40").Value = 3185
$ws.Range("K40").Value = 790
$ws.Range("L40").Value = 3185
$ws.Range("M40").Value = -654
$ws.Range("N40").Value = -3457

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 7988.6665
$ws.Range("J47").Value = 7988.6665
$ws.Range("L47").Value = 7988.6665
$ws.Range("N47").Value = -9132.666499999999

$ws.Range("H126").Value = 92981.91
$ws.Range("I126").Value = 127037.625
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 381112.875
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -378642.875
$ws.Range("N126").Value = -11440.0001
